$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 29 previously had no data (gap between row 28 and row 30).
# Fill it in with the new vessel entry - no shifting of existing rows.
$ws.Range("A29").Value = "Thunder"
$ws.Range("B29").Value = "Jackson Offshore"
$ws.Range("C29").Value = 252
$ws.Range("D29").Value = "OSV"

$ws.Range("A30").Select()
